# Swap the "group-code" and "group-name" columns (D and E) in the
# ReportingOrganisationGroup sheet, including the header row, so that the
# codeforiati:group-name column comes before codeforiati:group-code in the
# shared-string table ordering that Excel will produce, while keeping the
# actual data (code <-> name) correctly associated with its row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
